$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.880.19"
$ws.Range("E2").Value = "  -0.61%  "

$ws.Range("D3").Value = "2.750.27"
$ws.Range("E3").Value = "  -0.04%  "

$ws.Range("D4").Value = "'0.999"
$ws.Range("E4").Value = "  -0.14%  "

$ws.Range("D5").Value = "'579.07"
$ws.Range("E5").Value = "  -2.35%  "

$ws.Range("D6").Value = "'159.12"
$ws.Range("E6").Value = "  +4.10%  "

$ws.Range("E7").Value = "  +0.42%  "

$ws.Range("D8").Value = "'0.610"
$ws.Range("E8").Value = "  -0.44%  "

$ws.Range("D9").Value = "'0.112"
$ws.Range("E9").Value = "  -2.32%  "

$ws.Range("D10").Value = "'0.391"
$ws.Range("E10").Value = "  +0.23%  "

$ws.Range("B11").Value = "Toncoin"
$ws.Range("C11").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D11").Value = "'5.67"
$ws.Range("E11").Value = "  -16.56%  "

$ws.Range("B12").Value = "TRON"
$ws.Range("C12").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D12").Value = "'0.158"
$ws.Range("E12").Value = "  +0.07%  "

$ws.Range("D13").Value = "3.232.43"
$ws.Range("E13").Value = "  -0.16%  "

$ws.Range("D14").Value = "'27.01"
$ws.Range("E14").Value = "  +0.90%  "

$ws.Range("D15").Value = "63.760.65"
$ws.Range("E15").Value = "  -0.56%  "

$ws.Range("D16").Value = "'0.0000156"
$ws.Range("E16").Value = "  +1.42%  "

$ws.Range("D17").Value = "2.750.20"
$ws.Range("E17").Value = "  -0.94%  "

$ws.Range("D18").Value = "'12.31"
$ws.Range("E18").Value = "  +1.48%  "

$ws.Range("D19").Value = "'4.95"
$ws.Range("E19").Value = "  +0.53%  "

$ws.Range("D20").Value = "'361.08"
$ws.Range("E20").Value = "  -1.79%  "

$ws.Range("D21").Value = "'6.88"
$ws.Range("E21").Value = "  -2.60%  "

$ws.Range("D22").Value = "'0.567"
$ws.Range("E22").Value = "  +5.35%  "

$ws.Range("D23").Value = "'0.991"
$ws.Range("E23").Value = "  -0.35%  "

$ws.Range("D24").Value = "'66.34"
$ws.Range("E24").Value = "  -0.14%  "

$ws.Range("D25").Value = "'0.173"
$ws.Range("E25").Value = "  +1.73%  "

$ws.Range("D26").Value = "'8.68"
$ws.Range("E26").Value = "  -0.13%  "

$ws.Range("E27").Value = "  +0.21%  "

$ws.Range("D28").Value = "0.0₃0937"
$ws.Range("E28").Value = "  +2.00%  "

$ws.Range("D29").Value = "'2.01"
$ws.Range("E29").Value = "  -2.03%  "

$ws.Range("D30").Value = "'7.13"
$ws.Range("E30").Value = "  -0.28%  "

$ws.Range("E31").Value = "  +4.96%  "

$ws.Range("D32").Value = "'168.24"
$ws.Range("E32").Value = "  -2.36%  "

$ws.Range("D34").Value = "'5.03"
$ws.Range("E34").Value = "  +4.66%  "

$ws.Range("D35").Value = "'20.55"
$ws.Range("E35").Value = "  -0.95%  "

$ws.Range("D36").Value = "'1.47"
$ws.Range("E36").Value = "  +1.78%  "

$ws.Range("D37").Value = "'1.82"
$ws.Range("E37").Value = "  +0.54%  "

$ws.Range("D38").Value = "'1.02"
$ws.Range("E38").Value = "  -1.41%  "

$ws.Range("D39").Value = "'4.21"
$ws.Range("E39").Value = "  -1.31%  "

$ws.Range("D40").Value = "'6.15"
$ws.Range("E40").Value = "  +8.10%  "

$ws.Range("D41").Value = "'332.92"
$ws.Range("E41").Value = "  -5.28%  "

$ws.Range("D42").Value = "'39.70"
$ws.Range("E42").Value = "  +0.45%  "

$ws.Range("D43").Value = "'22.04"
$ws.Range("E43").Value = "  -1.51%  "

$ws.Range("D44").Value = "'0.0599"
$ws.Range("E44").Value = "  +0.82%  "

$ws.Range("D45").Value = "'21.99"
$ws.Range("E45").Value = "  -0.85%  "

$ws.Range("D46").Value = "'0.0259"
$ws.Range("E46").Value = "  -0.14%  "

$ws.Range("D47").Value = "'0.640"
$ws.Range("E47").Value = "  -2.20%  "

$ws.Range("D48").Value = "'137.15"
$ws.Range("E48").Value = "  -4.29%  "

$ws.Range("E49").Value = "  +0.09%  "

$ws.Range("E50").Value = "  +0.52%  "

$ws.Range("E51").Value = "  +1.09%  "
